$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Row 17 (Shell) : fill in time/measurement data ---
$ws.Range("G17").Value = 0.024305555555555556
$ws.Range("H17").Value = 0.3611111111111111
$ws.Range("I17").Value = 0.40277777777777773
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.006944444444444444
$ws.Range("M17").Value = 50

# --- Row 18 (Quick Sort) : fill in time/measurement data ---
$ws.Range("G18").Value = 0.024305555555555556
$ws.Range("H18").Value = 0.4166666666666667
$ws.Range("I18").Value = 0.46527777777777773
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.002777777777777778
$ws.Range("M18").Value = 50

# --- Row 24 (Ejecución de la Prueba) : fill in start/end times ---
$ws.Range("B24").Value = 0.041666666666666664
$ws.Range("C24").Value = 0.4479166666666667
$ws.Range("D24").Value = 0.513888888888889

# --- View state: scroll + selection moved ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M19").Select()

Write-Output "done"
